# Apply the "new version with timestamp" update described by the diff:
#  - Insert a new shortage item "NOVACTAM 1500MG VIAL I.M/I.V" as row 10 (item #4)
#  - Shift the two existing items (VIDROP..., سائل ريد) down by one row
#  - Recompute the totals row value (221 -> 267) and move the footer row down
#  - Bump the generated timestamp string by one minute (10:26 AM -> 10:27 AM)

function Set-TextValue($ws, $addr, $value) {
    # Force the cell to store the value as text (shared string) even when the
    # text looks like a plain number ("1", "0", "46.00", ...), while keeping
    # the cell's original number format / style untouched.
    $rng = $ws.Range($addr)
    $savedFormat = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = $savedFormat
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Footer row (was row 13) moves down to row 14 - copy formatting, values
#    and merged ranges.
# ---------------------------------------------------------------------------
$ws.Range("A13:Q13").Copy() | Out-Null
$ws.Range("A14:Q14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A13:Q13").Copy() | Out-Null
$ws.Range("A14:Q14").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Application.CutCopyMode = $false
$ws.Rows("14:14").RowHeight = 16.5
$ws.Range("A14:F14").Merge() | Out-Null
$ws.Range("G14:I14").Merge() | Out-Null
$ws.Range("K14:Q14").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 2) Totals row (was row 12) moves down to row 13, with the grand total
#    updated from 221 to 267 (the previous total plus the new item's price).
# ---------------------------------------------------------------------------
$ws.Range("A12:Q12").Copy() | Out-Null
$ws.Range("A13:Q13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Rows("13:13").RowHeight = 24.75
$ws.Range("P13").Value = 267
$ws.Range("Q13").Value = $null
$ws.Range("P13:Q13").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 3) Item rows: new row 12 becomes the old row 11 (سائل ريد / item 6).
# ---------------------------------------------------------------------------
$ws.Range("A11:Q11").Copy() | Out-Null
$ws.Range("A12:Q12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Rows("12:12").RowHeight = 25.5
$ws.Range("A12:B12").Merge() | Out-Null
$ws.Range("C12:G12").Merge() | Out-Null
$ws.Range("H12:K12").Merge() | Out-Null
$ws.Range("L12:M12").Merge() | Out-Null
$ws.Range("N12:O12").Merge() | Out-Null

Set-TextValue $ws "A12" 6
Set-TextValue $ws "C12" "سائل ريد"
Set-TextValue $ws "H12" "12:0"
Set-TextValue $ws "L12" "0"
Set-TextValue $ws "N12" "100.00"
Set-TextValue $ws "P12" "100.0000"
Set-TextValue $ws "Q12" "1:0"

# Row 11 becomes the old row 10 (VIDROP / item 5).
Set-TextValue $ws "C11" "VIDROP 2800 I.U./ML ORAL DROPS 15 ML"
Set-TextValue $ws "H11" "2:0"
Set-TextValue $ws "L11" "1"
Set-TextValue $ws "N11" "26.00"
Set-TextValue $ws "P11" "26.0000"

# Row 10 becomes the brand-new item (item 4).
Set-TextValue $ws "C10" "NOVACTAM 1500MG VIAL I.M/I.V"
Set-TextValue $ws "H10" "10:0"
Set-TextValue $ws "L10" "1"
Set-TextValue $ws "N10" "46.00"
Set-TextValue $ws "P10" "46.0000"

# ---------------------------------------------------------------------------
# 4) Bump the printed timestamp by one minute.
# ---------------------------------------------------------------------------
Set-TextValue $ws "A14" "Saturday, 26 July, 2025 10:27 AM"

Write-Host "Edit applied successfully"
